$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns L:N (old headers "COSTO", "PRECIO SIN IVA", "PRECIO CON IVA")
# This shifts old O (ITEM) -> L, old P (FECHA DESPACHO) -> M, old Q (blank) -> N
$ws.Range("L:N").Delete()

# The old "N° GUIA" header (now in K5) becomes "FIRMA"
$ws.Range("K5").Value = "FIRMA"

# Update the active selection to match the new layout
$ws.Range("L16").Select()
